# rtsp_cam.xlsx update (27 Aug) - fixed gpu problem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "gpu" flag column (C) for several rtsp cameras
$ws.Range("C2").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C9").Value = 0

# Move/restore the active selection to C4
$ws.Range("C4").Select()
